$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 16's formatting down to the new row 17 so the new row reuses the
# existing date-time style (A column) instead of minting a new numFmt/xf.
$ws.Range("A16:N16").Copy()
$ws.Range("A17:N17").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A17").Value = 42625.884247685186
$ws.Range("B17").Value = -32
$ws.Range("C17").Value = 61
$ws.Range("D17").Value = 38
$ws.Range("E17").Value = 28
$ws.Range("F17").Value = 71
$ws.Range("G17").Value = 11732
$ws.Range("H17").Value = 8861
$ws.Range("I17").Value = 426
$ws.Range("J17").Value = 125
$ws.Range("K17").Value = 78
$ws.Range("L17").Value = 2
$ws.Range("M17").Value = 5
$ws.Range("N17").Value = "Named"
